$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A1" = 0.171
    "B1" = -0.117
    "C1" = -0.123
    "D1" = -45
    "E1" = 10
    "F1" = -170
    "A2" = 0.147
    "B2" = 0.413
    "C2" = 0.184
    "D2" = 79
    "E2" = 30
    "F2" = 96
    "A3" = -0.351
    "B3" = 0.213
    "C3" = 0.1
    "D3" = -111
    "E3" = -74
    "F3" = -21
    "A4" = -0.304
    "B4" = -0.285
    "C4" = 0.002
    "D4" = 152
    "E4" = -52
    "F4" = -171
    "A5" = 0.134
    "B5" = -0.163
    "C5" = 0.8129999999999999
    "D5" = 64
    "E5" = 40
    "F5" = -6
    "A6" = -0.186
    "B6" = -0.615
    "C6" = 0.289
    "D6" = 171
    "E6" = 46
    "F6" = 0
    "A7" = -0.122
    "B7" = -0.055
    "C7" = 0.016
    "D7" = -141
    "E7" = -24
    "F7" = -59
    "A8" = -0.539
    "B8" = -0.038
    "C8" = 0.356
    "D8" = -63
    "E8" = 53
    "F8" = -176
    "A9" = 0.122
    "B9" = 0.43
    "C9" = 0.07099999999999999
    "D9" = 54
    "E9" = -31
    "F9" = 107
    "A10" = -0.316
    "B10" = -0.303
    "C10" = -0.226
    "D10" = -74
    "E10" = -8
    "F10" = -156
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
